# ==========================================================================
# Edit script: add "Player Info" sheet, rename MATCH_CARD_LINK columns to
# MATCH_CODE (storing bare match codes instead of full URLs), clean up a
# few stray empty cells, and append a new "ODI Batting Extra" sheet.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# Helper: apply the existing bold/bordered/centered header style (the
# same style already used by the header rows of "ODI Batting") to a
# target range, by copying format from a known header cell. This keeps
# re-using the same style index instead of creating new ones. The
# source range is re-fetched fresh every time (rather than cached in a
# variable) because once new sheets get inserted/moved, a previously
# captured Range reference can become stale.
# --------------------------------------------------------------------
function Apply-HeaderStyle($range) {
    $wb.Worksheets.Item("ODI Batting").Range("A1").Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# --------------------------------------------------------------------
# 1) Insert a new "Player Info" sheet before the existing first sheet.
# --------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"
Apply-HeaderStyle($playerInfo.Range("A1:D1"))

$playerInfo.Cells.Item(2,1).NumberFormat = "@"
$playerInfo.Cells.Item(2,1).Value = "3871"
$playerInfo.Cells.Item(2,2).Value = "Frank Dimuth Madushanka Karunaratne"
$playerInfo.Cells.Item(2,3).Value = "Left Handed"
$playerInfo.Cells.Item(2,4).Value = "Right Arm Medium"

# --------------------------------------------------------------------
# 2) "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE, store
#    bare match codes instead of full URLs, and clear a few stray empty
#    inning-number cells.
# --------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Cells.Item(1,4).Value = "MATCH_CODE"

$battingCodes = @{
    2 = "3305"; 3 = "3308"; 4 = "3576"; 5 = "3578"; 6 = "3590";
    7 = "3591"; 8 = "3723"; 9 = "3726"; 10 = "3729"; 11 = "3734";
    12 = "3737"; 13 = "3739"; 14 = "3743"; 15 = "3748"; 16 = "3759";
    17 = "3765"; 18 = "3769"; 19 = "4302"; 20 = "4305"; 21 = "4309";
    22 = "4322"; 23 = "4331"; 24 = "4339"; 25 = "4344"; 26 = "4350";
    27 = "4356"; 28 = "4357"; 29 = "4358"; 30 = "4413"; 31 = "4414";
    32 = "4417"; 33 = "4449"; 34 = "4450"; 35 = "4451"
}
$battingCodes.Keys | Sort-Object | ForEach-Object {
    $cell = $batting.Cells.Item($_, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$_]
}

# Remove the stray empty INNING_NUMBER cells on rows 10, 14, 17, 18.
$batting.Cells.Item(10,2).ClearContents()
$batting.Cells.Item(14,2).ClearContents()
$batting.Cells.Item(17,2).ClearContents()
$batting.Cells.Item(18,2).ClearContents()

# --------------------------------------------------------------------
# 3) "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and store
#    bare match codes instead of full URLs.
# --------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Cells.Item(1,2).Value = "MATCH_CODE"

$bowlingCodes = @{ 2 = "3590"; 3 = "4344" }
$bowlingCodes.Keys | Sort-Object | ForEach-Object {
    $cell = $bowling.Cells.Item($_, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$_]
}

# --------------------------------------------------------------------
# 4) Append a new "ODI Batting Extra" sheet at the end of the workbook.
# --------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extra.Cells.Item(1,1).Value = "MATCH_CODE"
$extra.Cells.Item(1,2).Value = "BATTING_POSITION"
$extra.Cells.Item(1,3).Value = "NUM_4"
$extra.Cells.Item(1,4).Value = "NUM_6"
$extra.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1,6).Value = "MAN_OF_MATCH"
Apply-HeaderStyle($extra.Range("A1:F1"))

# Row data: MatchCode, BattingPosition(number or $null), Num4, Num6, PercentRuns, ManOfMatch
$extraRows = @(
    @("3759", 4,     "3", "0", "9.75%",  "NO"),
    @("3765", 4,     $null, $null, $null, "NO"),
    @("3769", 6,     $null, $null, $null, "NO"),
    @("4302", $null, $null, $null, $null, "NO"),
    @("4305", 2,     "4", "0", "38.24%", "NO"),
    @("4309", 1,     "3", "0", "14.93%", "NO"),
    @("4322", 1,     "9", "0", "39.27%", "NO"),
    @("4331", 1,     "0", "0", "0.43%",  "NO"),
    @("4339", 1,     "0", "0", $null,    "NO"),
    @("4344", 1,     "4", "0", "9.47%",  "NO"),
    @("4350", 1,     "2", "0", "3.79%",  "NO"),
    @("4356", $null, $null, $null, $null, "NO"),
    @("4357", $null, $null, $null, $null, "NO"),
    @("4358", $null, $null, $null, $null, "NO"),
    @("4413", 2,     "7", "0", "17.93%", "NO"),
    @("4414", $null, $null, $null, $null, "NO"),
    @("4417", 2,     "5", "0", "14.33%", "NO"),
    @("4449", 2,     "4", "0", "22.41%", "NO"),
    @("4450", $null, $null, $null, $null, "NO"),
    @("4451", 2,     "1", "1", "11.31%", "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $codeCell = $extra.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    $posCell = $extra.Cells.Item($r, 2)
    if ($row[1] -ne $null) {
        $posCell.Value = $row[1]
    } else {
        $posCell.NumberFormat = "@"
    }

    $num4Cell = $extra.Cells.Item($r, 3)
    $num4Cell.NumberFormat = "@"
    if ($row[2] -ne $null) { $num4Cell.Value = $row[2] }

    $num6Cell = $extra.Cells.Item($r, 4)
    $num6Cell.NumberFormat = "@"
    if ($row[3] -ne $null) { $num6Cell.Value = $row[3] }

    $pctCell = $extra.Cells.Item($r, 5)
    $pctCell.NumberFormat = "@"
    if ($row[4] -ne $null) { $pctCell.Value = $row[4] }

    $momCell = $extra.Cells.Item($r, 6)
    $momCell.Value = $row[5]

    $r = $r + 1
}

Write-Host "Edit complete"
